# Update the "Время (мск)" time-slot labels in column C (rows 6-11).
# A new 5-minute slot ("12:25-12:30" / "12:30-12:35") was inserted, shifting
# the later "22:3x-22:5x" slots down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value  = "12:25-12:30"
$ws.Range("C7").Value  = "12:30-12:35"
$ws.Range("C8").Value  = "22:35-22:40"
$ws.Range("C9").Value  = "22:40-22:45"
$ws.Range("C10").Value = "22:45-22:50"
$ws.Range("C11").Value = "22:50-22:55"

# Move the active selection to C13, matching the saved view state.
$ws.Range("C13").Select()
